$wb = $excel.ActiveWorkbook

# ALC row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 7407728.5
$ws.Cells.Item(41, 9).Value = 9259573
$ws.Cells.Item(41, 10).Value = 350
$ws.Cells.Item(41, 11).Value = 9259573
$ws.Cells.Item(41, 12).Value = 350
$ws.Cells.Item(41, 13).Value = -9259133
$ws.Cells.Item(41, 14).Value = -1230

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 9430.093999999999
$ws.Cells.Item(62, 9).Value = 6446.2173
$ws.Cells.Item(62, 11).Value = 6446.2173
$ws.Cells.Item(62, 13).Value = -5822.2173

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 9430.093999999999
$ws.Cells.Item(65, 9).Value = 6446.2173
$ws.Cells.Item(65, 11).Value = 32231.0865
$ws.Cells.Item(65, 13).Value = -29111.0865

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 1139.3636
$ws.Cells.Item(86, 9).Value = 661.5
$ws.Cells.Item(86, 10).Value = 1712.8
$ws.Cells.Item(86, 11).Value = 661.5
$ws.Cells.Item(86, 12).Value = 1712.8
$ws.Cells.Item(86, 13).Value = 461.5
$ws.Cells.Item(86, 14).Value = -3958.8

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 1139.3636
$ws.Cells.Item(89, 9).Value = 661.5
$ws.Cells.Item(89, 10).Value = 1712.8
$ws.Cells.Item(89, 11).Value = 3307.5
$ws.Cells.Item(89, 12).Value = 8564
$ws.Cells.Item(89, 13).Value = 2308.5
$ws.Cells.Item(89, 14).Value = -19796

# ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 585446.9399999999
$ws.Cells.Item(92, 9).Value = 741232.4399999999
$ws.Cells.Item(92, 10).Value = 1251.25
$ws.Cells.Item(92, 11).Value = 741232.4399999999
$ws.Cells.Item(92, 12).Value = 1251.25
$ws.Cells.Item(92, 13).Value = -739984.4399999999
$ws.Cells.Item(92, 14).Value = -3747.25

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 4276535.5
$ws.Cells.Item(106, 9).Value = 4833879.5
$ws.Cells.Item(106, 11).Value = 4833879.5
$ws.Cells.Item(106, 13).Value = -4833248.5

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 427707.84
$ws.Cells.Item(107, 9).Value = 463305.6
$ws.Cells.Item(107, 10).Value = 535
$ws.Cells.Item(107, 11).Value = 463305.6
$ws.Cells.Item(107, 12).Value = 535
$ws.Cells.Item(107, 13).Value = -461385.6
$ws.Cells.Item(107, 14).Value = -4375

# ARM row 4
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(4, 14).ClearContents()

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1022.6667
$ws.Cells.Item(45, 9).Value = 1041.4286
$ws.Cells.Item(45, 10).Value = 996.4
$ws.Cells.Item(45, 11).Value = 1041.4286
$ws.Cells.Item(45, 12).Value = 996.4
$ws.Cells.Item(45, 13).Value = -664.4286
$ws.Cells.Item(45, 14).Value = -1750.4

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2989.8667
$ws.Cells.Item(20, 9).Value = 3234.8
$ws.Cells.Item(20, 10).Value = 2500
$ws.Cells.Item(20, 11).Value = 3234.8
$ws.Cells.Item(20, 12).Value = 2500
$ws.Cells.Item(20, 13).Value = -2987.8
$ws.Cells.Item(20, 14).Value = -2994

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2269.0527
$ws.Cells.Item(134, 9).Value = 1397.5238
$ws.Cells.Item(134, 10).Value = 4709.3335
$ws.Cells.Item(134, 11).Value = 4192.5714
$ws.Cells.Item(134, 12).Value = 14128.0005
$ws.Cells.Item(134, 13).Value = -1657.5714
$ws.Cells.Item(134, 14).Value = -19198.0005

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 533.7143
$ws.Cells.Item(22, 9).Value = 477.83334
$ws.Cells.Item(22, 10).Value = 608.2222
$ws.Cells.Item(22, 11).Value = 477.83334
$ws.Cells.Item(22, 12).Value = 608.2222
$ws.Cells.Item(22, 13).Value = -127.83334
$ws.Cells.Item(22, 14).Value = -1308.2222

# CRP row 25
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(25, 8).Value = 8000
$ws.Cells.Item(25, 10).Value = 8000
$ws.Cells.Item(25, 12).Value = 8000
$ws.Cells.Item(25, 14).Value = -8348

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3900.4443
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 3900.4443
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 12).Value = 3900.4443
$ws.Cells.Item(31, 13).ClearContents()
$ws.Cells.Item(31, 14).Value = -4490.4443

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3900.4443
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 3900.4443
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 3900.4443
$ws.Cells.Item(34, 13).ClearContents()
$ws.Cells.Item(34, 14).Value = -4304.4443

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 10417800
$ws.Cells.Item(99, 9).Value = 15625950
$ws.Cells.Item(99, 11).Value = 15625950
$ws.Cells.Item(99, 13).Value = -15624452

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 2897
$ws.Cells.Item(122, 9).Value = 1334.4
$ws.Cells.Item(122, 11).Value = 4003.2
$ws.Cells.Item(122, 13).Value = -1553.2

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 10417800
$ws.Cells.Item(126, 9).Value = 15625950
$ws.Cells.Item(126, 11).Value = 46877850
$ws.Cells.Item(126, 13).Value = -46875380

# CRP row 133
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(133, 8).Value = 26711.758
$ws.Cells.Item(133, 10).Value = 28631.385
$ws.Cells.Item(133, 12).Value = 28631.385
$ws.Cells.Item(133, 14).Value = -33691.38499999999

# CRP row 135
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(135, 8).Value = 39887.93
$ws.Cells.Item(135, 10).Value = 39887.93
$ws.Cells.Item(135, 12).Value = 39887.93
$ws.Cells.Item(135, 14).Value = -50027.93

# CUL row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 1200
$ws.Cells.Item(80, 10).Value = 1200
$ws.Cells.Item(80, 12).Value = 3600
$ws.Cells.Item(80, 14).Value = -5472

# CUL row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(83, 8).Value = 1200
$ws.Cells.Item(83, 10).Value = 1200
$ws.Cells.Item(83, 12).Value = 10800
$ws.Cells.Item(83, 14).Value = -20160

# GSM row 11
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 13116616
$ws.Cells.Item(11, 10).Value = 5000
$ws.Cells.Item(11, 12).Value = 5000
$ws.Cells.Item(11, 14).Value = -5278

# GSM row 40
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 14).ClearContents()

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 11560.77
$ws.Cells.Item(61, 9).Value = 12389.05
$ws.Cells.Item(61, 10).Value = 8799.833000000001
$ws.Cells.Item(61, 11).Value = 12389.05
$ws.Cells.Item(61, 12).Value = 8799.833000000001
$ws.Cells.Item(61, 13).Value = -12187.05
$ws.Cells.Item(61, 14).Value = -9203.833000000001

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 11560.77
$ws.Cells.Item(113, 9).Value = 12389.05
$ws.Cells.Item(113, 10).Value = 8799.833000000001
$ws.Cells.Item(113, 11).Value = 12389.05
$ws.Cells.Item(113, 12).Value = 8799.833000000001
$ws.Cells.Item(113, 13).Value = -10219.05
$ws.Cells.Item(113, 14).Value = -13139.833

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 3756.549
$ws.Cells.Item(136, 9).Value = 2300.375
$ws.Cells.Item(136, 10).Value = 9051.727999999999
$ws.Cells.Item(136, 11).Value = 6901.125
$ws.Cells.Item(136, 12).Value = 27155.184
$ws.Cells.Item(136, 13).Value = -4351.125
$ws.Cells.Item(136, 14).Value = -32255.184

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 17816.373
$ws.Cells.Item(136, 9).Value = 19312.889
$ws.Cells.Item(136, 10).Value = 1654
$ws.Cells.Item(136, 11).Value = 57938.667
$ws.Cells.Item(136, 12).Value = 4962
$ws.Cells.Item(136, 13).Value = -55388.667
$ws.Cells.Item(136, 14).Value = -10062
